$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of an existing wrapped-text cell (D4, style index 6)
# onto the new rows 5:6 range B:F, so we reuse the existing cellXf instead of
# creating a brand-new one.
$ws.Range("D4").Copy()
$ws.Range("B5:F6").PasteSpecial(-4122)  # xlPasteFormats

# Match the row heights used by the author for the two new rows.
$ws.Rows.Item(5).RowHeight = 28.8
$ws.Rows.Item(6).RowHeight = 28.8

# Row 5: Aggressive Cows - Linear Search
$ws.Range("D5").Value = "Aggressive Cows - Linear Search"
# Row 6: Aggressive Cows - Binary Search
$ws.Range("D6").Value = "Aggressive Cows - Binary Search"
# Question-number column for both new rows
$ws.Range("C5").Value = "Binary S2 4"
$ws.Range("C6").Value = "Binary S2 4"
# S.no. values
$ws.Range("B5").Value = 3
$ws.Range("B6").Value = 4
# Page no. in notes column
$ws.Range("E5").Value = "na"
$ws.Range("E6").Value = "na"

# Match final selection left behind by the author
$ws.Range("C7").Select()
